# hydraulic gear system - add "Nose gear door" and "Nose gear" actuator
# blocks to the Actuators sheet, add "Max force retract/extend" columns
# to the existing blocks, tweak the "Main gear door" bore/rod diameter,
# add two new threaded comments and re-fill the bore/rod-side cells that
# were highlighted red to match the orange highlight used by the other
# coordinate cells in the same rows. Also nudge the landing-gear-door
# clearance diagram further right/down to make room for the new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Actuators")

# ---------------------------------------------------------------------
# 1. Existing "Main gear door" actuator block (rows 43-47): update the
#    rod diameter formula/value and add the "Max force retract/extend"
#    helper columns that now also appear on the other three blocks.
# ---------------------------------------------------------------------
$ws.Range("C44").Value = 0.0538
$ws.Range("D44").Formula = "=0.03015"

$ws.Range("H32").Value = "Max force extend"
$ws.Range("I32").Value = "Max force retract"
$ws.Range("H33").Formula = "=(E30-F30)*20684000"
$ws.Range("I33").Formula = "=E30*20684000"

$ws.Range("H39").Value = "Max force retract"
$ws.Range("I39").Value = "Max force extend"
$ws.Range("H40").Formula = "=(E37-F37)*20684000"
$ws.Range("I40").Formula = "=E37*20684000"

$ws.Range("H46").Value = "Max force retract"
$ws.Range("I46").Value = "Max force extend"
$ws.Range("H47").Formula = "=(E44-F44)*20684000"
$ws.Range("I47").Formula = "=E44*20684000"

# ---------------------------------------------------------------------
# 2. New "Nose gear door" actuator block (rows 50-56), formatted the
#    same way as the "Main gear door" block (rows 43-47).
# ---------------------------------------------------------------------
$ws.Range("A43:J47").Copy()
$ws.Range("A50:J54").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("B50").Value = "Displacement length (m)"
$ws.Range("C50").Value = "piston diameter (m)"
$ws.Range("D50").Value = "rod diameter (m)"
$ws.Range("E50").Value = "Piston area (m2)"
$ws.Range("F50").Value = "Rod area (m2)"
$ws.Range("G50").Value = "Volume bore side (m3)"
$ws.Range("H50").Value = "Volume rod side (m3)"
$ws.Range("I50").Value = "Volume Ratio extended"
$ws.Range("J50").Value = "Volume Impact On Reservoir (L)"

$ws.Range("A51").Value = "Nose gear door"
$ws.Range("B51").Value = 0.162
$ws.Range("C51").Value = 0.0378
$ws.Range("D51").Formula = "=0.023"
$ws.Range("E51").Formula = "=PI()* (C51/2)^2"
$ws.Range("F51").Formula = "=PI()* (D51/2)^2"
$ws.Range("G51").Formula = "=E51*B51"
$ws.Range("H51").Formula = "=(E51-F51)*B51"
$ws.Range("I51").Formula = "=G51/H51"
$ws.Range("J51").Formula = "=(G51-H51) * 1000"

$ws.Range("A52").Value = "Both actuators"
$ws.Range("G52").Formula = "=G51*2"
$ws.Range("H52").Formula = "=H51*2"
$ws.Range("I52").Formula = "=G52/H52"
$ws.Range("J52").Formula = "=(G52-H52) * 1000"

$ws.Range("B53").Value = "Actuator Length(m)"
$ws.Range("C53").Value = "Control Arm position X"
$ws.Range("D53").Value = "Control Arm position Y"
$ws.Range("E53").Value = "Anchor point position X"
$ws.Range("F53").Value = "Anchor point position Y"
$ws.Range("H53").Value = "Max force retract"
$ws.Range("I53").Value = "Max force extend"

$ws.Range("C54").Value = -0.1465
$ws.Range("D54").Value = 0
$ws.Range("E54").Value = -0.1465
$ws.Range("F54").Value = 0.4
$ws.Range("H54").Formula = "=(E51-F51)*20684000"
$ws.Range("I54").Formula = "=E51*20684000"

$ws.Range("C55").Value = 'Note those coordinates are chosen "randomly" so we get the correct actuator travel of 162mm'
$ws.Range("C56").Value = "This actuator in fact should push a crankbell that gives a particular motion ratio."

# New threaded comment on D54, duplicated from the one on D47.
$ws.Range("D54").AddCommentThreaded("This value gives correct actuator displacement") | Out-Null

# ---------------------------------------------------------------------
# 3. New "Nose gear" actuator block (rows 58-63), same layout again.
# ---------------------------------------------------------------------
$ws.Range("A43:J47").Copy()
$ws.Range("A58:J62").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("B58").Value = "Displacement length (m)"
$ws.Range("C58").Value = "piston diameter (m)"
$ws.Range("D58").Value = "rod diameter (m)"
$ws.Range("E58").Value = "Piston area (m2)"
$ws.Range("F58").Value = "Rod area (m2)"
$ws.Range("G58").Value = "Volume bore side (m3)"
$ws.Range("H58").Value = "Volume rod side (m3)"
$ws.Range("I58").Value = "Volume Ratio extended"
$ws.Range("J58").Value = "Volume Impact On Reservoir (L)"

$ws.Range("A59").Value = "Nose gear"
$ws.Range("B59").Value = 0.32
$ws.Range("C59").Value = 0.0792
$ws.Range("D59").Formula = "=0.035"
$ws.Range("E59").Formula = "=PI()* (C59/2)^2"
$ws.Range("F59").Formula = "=PI()* (D59/2)^2"
$ws.Range("G59").Formula = "=E59*B59"
$ws.Range("H59").Formula = "=(E59-F59)*B59"
$ws.Range("I59").Formula = "=G59/H59"
$ws.Range("J59").Formula = "=(G59-H59) * 1000"

$ws.Range("A60").Value = "Both actuators"
$ws.Range("G60").Formula = "=G59*2"
$ws.Range("H60").Formula = "=H59*2"
$ws.Range("I60").Formula = "=G60/H60"
$ws.Range("J60").Formula = "=(G60-H60) * 1000"

$ws.Range("B61").Value = "Actuator Length(m)"
$ws.Range("C61").Value = "Control Arm position Z"
$ws.Range("D61").Value = "Control Arm position Y"
$ws.Range("E61").Value = "Anchor point position Z"
$ws.Range("F61").Value = "Anchor point position Y"
$ws.Range("H61").Value = "Max force retract"
$ws.Range("I61").Value = "Max force extend"

$ws.Range("C62").Value = 0.212
$ws.Range("D62").Value = -0.093
$ws.Range("E62").Value = 0
$ws.Range("F62").Value = 0.56
$ws.Range("H62").Formula = "=(E59-F59)*20684000"
$ws.Range("I62").Formula = "=E59*20684000"

$ws.Range("C63").Value = 'Note those coordinates are chosen "randomly" so we get the correct actuator travel of 320mm'

# New threaded comment on D62, duplicated from the one on D47.
$ws.Range("D62").AddCommentThreaded("This value gives correct actuator displacement") | Out-Null

# ---------------------------------------------------------------------
# 4. Re-colour the bore/rod-side coordinate cells (columns C & E) in the
#    three original blocks so they match the orange highlight already
#    used on columns D & F in those same rows.
# ---------------------------------------------------------------------
$ws.Range("C33").Interior.Color = $ws.Range("D33").Interior.Color
$ws.Range("E33").Interior.Color = $ws.Range("D33").Interior.Color
$ws.Range("C40").Interior.Color = $ws.Range("D40").Interior.Color
$ws.Range("E40").Interior.Color = $ws.Range("D40").Interior.Color
$ws.Range("C47").Interior.Color = $ws.Range("D47").Interior.Color
$ws.Range("E47").Interior.Color = $ws.Range("D47").Interior.Color

# ---------------------------------------------------------------------
# 5. Move the landing-gear-door clearance picture further right/down so
#    it doesn't overlap the newly added rows.
# ---------------------------------------------------------------------
$shp = $ws.Shapes.Item(3)
$shp.Left = $shp.Left + 5692589 / 9525
$shp.Top = $shp.Top - 1199030 / 9525

# ---------------------------------------------------------------------
# 6. Misc view/selection bookkeeping to mirror the authored session.
# ---------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 37
$ws.Range("G62").Select() | Out-Null
